$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Name" column)
$ws.Columns("B:B").Insert()

# The inserted column did not inherit the header formatting (bold, border,
# centered) from its neighbour, so copy that formatting explicitly before
# writing the new header text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "Name"

# The inserted column copied column A's formatting (border/style) for the
# data rows; the target layout keeps those cells unstyled (like the other
# plain data columns), so clear the inherited formatting first.
$ws.Range("B2:B6").ClearFormats()

# Fill the new "Name" column with the same values as column A for each row
$ws.Range("B2").Value = "idle1"
$ws.Range("B3").Value = "drive1"
$ws.Range("B4").Value = "live1"
$ws.Range("B5").Value = "ufpe1"
$ws.Range("B6").Value = "long1"

# Update the recalculated "max Leistung" (H) and "max Drehmoment" (I) values
$ws.Range("H3").Value = 62.11
$ws.Range("I3").Value = 296.48

$ws.Range("H4").Value = 70.22
$ws.Range("I4").Value = 389.5

$ws.Range("H5").Value = 54.04
$ws.Range("I5").Value = 328.97

$ws.Range("H6").Value = 99.84999999999999
$ws.Range("I6").Value = 451.39
